# Actualización automática 2025-11-25 17:30:09
#
# Updates sales figures across the three worksheets of the workbook:
#   - "VENTAS POR GRUPO"      (per-client sales broken down by product group)
#   - "VENTA MENSUAL"         (per-client sales broken down by month)
#   - "CUMPLIMIENTO MENSUAL"  (budget-vs-sales compliance summary by group)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("L9").Value  = 980.78
$wsGrupo.Range("D12").Value = 1425.6
$wsGrupo.Range("H16").Value = 670.91
$wsGrupo.Range("I24").Value = 259.2
$wsGrupo.Range("M24").Value = 190.32
$wsGrupo.Range("D52").Value = 475.2
$wsGrupo.Range("D53").Value = 633.6
$wsGrupo.Range("D56").Value = "10 de 54"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F9").Value  = 924.9400000000001
$wsMensual.Range("F12").Value = 4643.83
$wsMensual.Range("F16").Value = 8525.76
$wsMensual.Range("F24").Value = 4830.68
$wsMensual.Range("F53").Value = 2325.71
$wsMensual.Range("F54").Value = 2325.71
$wsMensual.Range("F55").Value = 657
$wsMensual.Range("F56").Value = 657
$wsMensual.Range("F60").Value = 80574.56999999999

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value  = 8962.85
$wsCumplimiento.Range("E3").Value  = -2339.59
$wsCumplimiento.Range("F3").Value  = 1.353238435453236

$wsCumplimiento.Range("D6").Value  = 3979.31
$wsCumplimiento.Range("E6").Value  = -1071.72631853974
$wsCumplimiento.Range("F6").Value  = 1.368596895550567

$wsCumplimiento.Range("D7").Value  = 727.2
$wsCumplimiento.Range("E7").Value  = 592.8
$wsCumplimiento.Range("F7").Value  = 0.5509090909090909

$wsCumplimiento.Range("D11").Value = 19825.22
$wsCumplimiento.Range("E11").Value = -5589.230000000001
$wsCumplimiento.Range("F11").Value = 1.392612666909713

$wsCumplimiento.Range("D12").Value = 41520.44
$wsCumplimiento.Range("E12").Value = 23423.56
$wsCumplimiento.Range("F12").Value = 0.6393268046316827

$wsCumplimiento.Range("D14").Value = 77591.86
$wsCumplimiento.Range("E14").Value = 21364.39685923838
$wsCumplimiento.Range("F14").Value = 0.7841026172844388
